{"js": "// Update the date line and the twenty-five multiplication problems in the\n// two-digit x two-digit practice sheet. Each old value is unique in the\n// document, so a simple ordered search/replace is safe and unambiguous.\nconst replacements = [\n  [\"2024-04-10 Wednesday\", \"2024-04-11 Thursday\"],\n  [\"58\u00d745=\", \"36\u00d719=\"],\n  [\"69\u00d783=\", \"32\u00d719=\"],\n  [\"24\u00d756=\", \"47\u00d758=\"],\n  [\"81\u00d771=\", \"77\u00d741=\"],\n  [\"98\u00d793=\", \"37\u00d733=\"],\n  [\"28\u00d788=\", \"12\u00d798=\"],\n  [\"52\u00d758=\", \"40\u00d735=\"],\n  [\"61\u00d794=\", \"56\u00d772=\"],\n  [\"19\u00d729=\", \"75\u00d740=\"],\n  [\"44\u00d783=\", \"46\u00d782=\"],\n  [\"34\u00d744=\", \"48\u00d761=\"],\n  [\"48\u00d786=\", \"52\u00d756=\"],\n  [\"11\u00d775=\", \"69\u00d767=\"],\n  [\"15\u00d763=\", \"48\u00d715=\"],\n  [\"51\u00d743=\", \"66\u00d775=\"],\n  [\"47\u00d734=\", \"59\u00d782=\"],\n  [\"91\u00d778=\", \"79\u00d787=\"],\n  [\"26\u00d748=\", \"93\u00d782=\"],\n  [\"13\u00d733=\", \"90\u00d748=\"],\n  [\"67\u00d754=\", \"91\u00d775=\"],\n  [\"49\u00d783=\", \"36\u00d739=\"],\n  [\"36\u00d796=\", \"84\u00d775=\"],\n  [\"14\u00d715=\", \"60\u00d746=\"],\n  [\"81\u00d739=\", \"71\u00d791=\"],\n  [\"46\u00d775=\", \"18\u00d746=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the twenty-five multiplication problems in the\n# two-digit x two-digit practice sheet. Each old value is unique in the\n# document, so a simple ordered Find/Replace is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-04-10 Wednesday\", \"2024-04-11 Thursday\"),\n  @(\"58\u00d745=\", \"36\u00d719=\"),\n  @(\"69\u00d783=\", \"32\u00d719=\"),\n  @(\"24\u00d756=\", \"47\u00d758=\"),\n  @(\"81\u00d771=\", \"77\u00d741=\"),\n  @(\"98\u00d793=\", \"37\u00d733=\"),\n  @(\"28\u00d788=\", \"12\u00d798=\"),\n  @(\"52\u00d758=\", \"40\u00d735=\"),\n  @(\"61\u00d794=\", \"56\u00d772=\"),\n  @(\"19\u00d729=\", \"75\u00d740=\"),\n  @(\"44\u00d783=\", \"46\u00d782=\"),\n  @(\"34\u00d744=\", \"48\u00d761=\"),\n  @(\"48\u00d786=\", \"52\u00d756=\"),\n  @(\"11\u00d775=\", \"69\u00d767=\"),\n  @(\"15\u00d763=\", \"48\u00d715=\"),\n  @(\"51\u00d743=\", \"66\u00d775=\"),\n  @(\"47\u00d734=\", \"59\u00d782=\"),\n  @(\"91\u00d778=\", \"79\u00d787=\"),\n  @(\"26\u00d748=\", \"93\u00d782=\"),\n  @(\"13\u00d733=\", \"90\u00d748=\"),\n  @(\"67\u00d754=\", \"91\u00d775=\"),\n  @(\"49\u00d783=\", \"36\u00d739=\"),\n  @(\"36\u00d796=\", \"84\u00d775=\"),\n  @(\"14\u00d715=\", \"60\u00d746=\"),\n  @(\"81\u00d739=\", \"71\u00d791=\"),\n  @(\"46\u00d775=\", \"18\u00d746=\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
